# "support for error cells"
#
# Adds a new worksheet named "error" at the end of the workbook
# (after "float" and "int") that demonstrates cells holding error
# values: #VALUE!, #NAME? and #DIV/0!, produced by real formulas.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() inserts at the front by default, so add it and then
# move it to be the last tab in the workbook.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "error"
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch the sheet by name (the reference returned by Add/Move can go
# stale once the sheet collection is reordered).
$ws = $wb.Worksheets.Item("error")

# Header / label cell (stored as a shared string).
$ws.Range("A1").Value = "errors"

# A2: text + number -> #VALUE!
$ws.Range("A2").Formula = "=A1+1"

# A3: calling an unknown name -> #NAME?
$ws.Range("A3").Formula = "=undefined_function"

# A4: division by zero -> #DIV/0!
$ws.Range("A4").Formula = "=10/0"

# Leave the selection cursor on the cell right below the data (A5), then
# restore "float" as the active/selected sheet/tab, matching the
# original workbook's active sheet.
$ws.Range("A5").Select()
$wb.Worksheets.Item("float").Activate()
